# Update LR-pair data rows (rows 2-5 changed, rows 6-9 added) per Dr Hou's advice
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = New-Object "object[]" 8
$rows[0] = @("FAPs", "Bdnf", "Ntrk2", "ECs", 3, 1, 2.020961333333334, 6.062884, 0.3447258214530571, 0.3447258214530571, 2, 0.6666666666666666, 0.1293823333333333, 0.388147, 0.006916257066299044, 0.006916257066299042, 0.2614766928831112, 2.353290235948, 0.002384212398560449, 0.002384212398560448)
$rows[1] = @("FAPs", "Bdnf", "Ntrk2", "FAPs", 3, 1, 2.020961333333334, 6.062884, 0.3447258214530571, 0.3447258214530571, 3, 1, 14.52590566666666, 43.57771699999999, 0.7764962582074056, 0.7764962582074055, 29.35629368398089, 264.206643155828, 0.267678310465773, 0.267678310465773)
$rows[2] = @("FAPs", "Bdnf", "Ntrk2", "M2", 3, 1, 2.020961333333334, 6.062884, 0.3447258214530571, 0.3447258214530571, 1, 0.3333333333333333, 0.01202633333333333, 0.036079, 0.0006428792150783161, 0.000642879215078316, 0.02430475464844445, 0.218742791836, 0.0002216170655129691, 0.000221617065512969)
$rows[3] = @("FAPs", "Bdnf", "Ntrk2", "sCs", 3, 1, 2.020961333333334, 6.062884, 0.3447258214530571, 0.3447258214530571, 3, 1, 4.039673, 12.119019, 0.215944605511217, 0.215944605511217, 8.164022932310667, 73.47620639079601, 0.07444168152321065, 0.07444168152321064)
$rows[4] = @("sCs", "Bdnf", "Ntrk2", "ECs", 3, 1, 3.841556666666667, 11.52467, 0.6552741785469429, 0.6552741785469429, 2, 0.6666666666666666, 0.1293823333333333, 0.388147, 0.006916257066299044, 0.006916257066299042, 0.4970295651655556, 4.473266086490001, 0.004532044667738595, 0.004532044667738594)
$rows[5] = @("sCs", "Bdnf", "Ntrk2", "FAPs", 3, 1, 3.841556666666667, 11.52467, 0.6552741785469429, 0.6552741785469429, 3, 1, 14.52590566666666, 43.57771699999999, 0.7764962582074056, 0.7764962582074055, 55.80208975315444, 502.2188077783899, 0.5088179477416326, 0.5088179477416326)
$rows[6] = @("sCs", "Bdnf", "Ntrk2", "M2", 3, 1, 3.841556666666667, 11.52467, 0.6552741785469429, 0.6552741785469429, 1, 0.3333333333333333, 0.01202633333333333, 0.036079, 0.0006428792150783161, 0.000642879215078316, 0.04619984099222223, 0.41579856893, 0.000421262149565347, 0.000421262149565347)
$rows[7] = @("sCs", "Bdnf", "Ntrk2", "sCs", 3, 1, 3.841556666666667, 11.52467, 0.6552741785469429, 0.6552741785469429, 3, 1, 4.039673, 12.119019, 0.215944605511217, 0.215944605511217, 15.51863274430333, 139.66769469873, 0.1415029239880064, 0.1415029239880064)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")
$startRow = 2

for ($ri = 0; $ri -lt $rows.Length; $ri++) {
    $rowVals = $rows[$ri]
    $r = $startRow + $ri
    for ($ci = 0; $ci -lt $cols.Length; $ci++) {
        $addr = $cols[$ci] + $r
        $ws.Range($addr).Value2 = $rowVals[$ci]
    }
}

Write-Output "Done updating rows 2-9"
